# Updated symbol list on Thu Dec 22 18:41:14 UTC 2022 with GitHub Actions
#
# Applies the latest price/volume-label refresh to the cryptos sheet.
# Price cells in column D are stored as plain text (not numbers), so each
# write temporarily forces a Text number format to stop Excel from
# re-interpreting the numeric-looking string as a float, then restores the
# cell's original "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $ref, $val) {
    $cell = $sheet.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2"  "241.85"
Set-TextValue $ws "D3"  "21.80"
Set-TextValue $ws "D4"  "5.396"
Set-TextValue $ws "D5"  "0.05678"
Set-TextValue $ws "D7"  "6.279"
Set-TextValue $ws "D8"  "0.8066"
Set-TextValue $ws "D9"  "0.8685"
Set-TextValue $ws "D10" "0.1425"
Set-TextValue $ws "D11" "0.07283"
Set-TextValue $ws "D12" "0.03061"
Set-TextValue $ws "D14" "0.09346"
Set-TextValue $ws "D15" "3.905"
Set-TextValue $ws "D16" "0.001589"
Set-TextValue $ws "D17" "0.04816"
Set-TextValue $ws "D18" "0.0005825"
Set-TextValue $ws "D19" "0.006334"
Set-TextValue $ws "D20" "0.0009989"
Set-TextValue $ws "D22" "0.0001500"
Set-TextValue $ws "D23" "3.734"
Set-TextValue $ws "D24" "2.150"
Set-TextValue $ws "D25" "0.3259"
Set-TextValue $ws "D26" "0.1311"
Set-TextValue $ws "D27" "0.0004003"
Set-TextValue $ws "D40" "0.03800"
Set-TextValue $ws "D41" "0.006687"
$ws.Range("E41").Value = "40KickTokenKICK"
Set-TextValue $ws "D42" "0.1047"
Set-TextValue $ws "D43" "0.002682"
Set-TextValue $ws "D44" "0.006829"
Set-TextValue $ws "D45" "0.00005609"
Set-TextValue $ws "D47" "0.5806"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"

Write-Output "Applied cryptos symbol list refresh"
